$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(756).Insert()

$ws.Cells.Item(756, 1).Value = 3
$ws.Cells.Item(756, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(756, 3).Value = "Coquimbo"
$ws.Cells.Item(756, 4).Value = 44939
$ws.Cells.Item(756, 5).Value = 5
$ws.Cells.Item(756, 6).Value = 100112045
$ws.Cells.Item(756, 7).Value = "Zapallo"
$ws.Cells.Item(756, 8).Value = "Camote"
$ws.Cells.Item(756, 9).Value = "1a nueva(o)"
$ws.Cells.Item(756, 10).Value = 175
$ws.Cells.Item(756, 11).Value = 650
$ws.Cells.Item(756, 12).Value = 700
$ws.Cells.Item(756, 13).Value = 677
$ws.Cells.Item(756, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(756, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(756, 16).Value = 677
$ws.Cells.Item(756, 17).Value = 1
$ws.Cells.Item(756, 18).Value = "Hortaliza"
